# Add data for 2022-10-08
# Updates year-to-date crime counts across the Citywide Totals, By Neighborhood,
# and each per-neighborhood sheet to reflect one additional day of data.

$wb = $excel.ActiveWorkbook

# Sheet 1: Citywide Totals
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("H2").Value = 91
$ws.Range("G3").Value = 111
$ws.Range("H3").Value = 114
$ws.Range("B9").Value = 314
$ws.Range("C9").Value = 389
$ws.Range("D9").Value = 335
$ws.Range("E9").Value = 355
$ws.Range("F9").Value = 425
$ws.Range("B10").Value = 1061
$ws.Range("C10").Value = 1275
$ws.Range("D10").Value = 1441
$ws.Range("E10").Value = 1757
$ws.Range("F10").Value = 1782
$ws.Range("G10").Value = 808
$ws.Range("H10").Value = 469
$ws.Range("I10").Value = 709
$ws.Range("B11").Value = 1487
$ws.Range("C11").Value = 1800
$ws.Range("D11").Value = 1969
$ws.Range("E11").Value = 2294
$ws.Range("F11").Value = 2385
$ws.Range("G11").Value = 1385
$ws.Range("H11").Value = 1064
$ws.Range("I11").Value = 1421

# Sheet 2: By Neighborhood
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("F5").Value = 31
$ws.Range("I5").Value = 27
$ws.Range("H9").Value = 3
$ws.Range("B11").Value = 9
$ws.Range("D19").Value = 50
$ws.Range("G19").Value = 35
$ws.Range("F28").Value = 94
$ws.Range("D32").Value = 84
$ws.Range("E32").Value = 111
$ws.Range("G32").Value = 92
$ws.Range("C36").Value = 82
$ws.Range("C47").Value = 61
$ws.Range("H47").Value = 37
$ws.Range("F50").Value = 53
$ws.Range("C53").Value = 292
$ws.Range("E53").Value = 575
$ws.Range("F53").Value = 542
$ws.Range("H61").Value = 17
$ws.Range("H62").Value = 12
$ws.Range("I63").Value = 5
$ws.Range("H74").Value = 19
$ws.Range("F76").Value = 54
$ws.Range("B78").Value = 34
$ws.Range("E90").Value = 9
$ws.Range("E93").Value = 5
$ws.Range("D97").Value = 21
$ws.Range("B99").Value = 1487
$ws.Range("C99").Value = 1800
$ws.Range("D99").Value = 1969
$ws.Range("E99").Value = 2294
$ws.Range("F99").Value = 2385
$ws.Range("G99").Value = 1385
$ws.Range("H99").Value = 1064
$ws.Range("I99").Value = 1421

# Sheet 3: Rogers Park
$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("F7").Value = 39
$ws.Range("F8").Value = 54

# Sheet 12: Garfield Park
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("E7").Value = 36
$ws.Range("D8").Value = 41
$ws.Range("G8").Value = 49
$ws.Range("D9").Value = 84
$ws.Range("E9").Value = 111
$ws.Range("G9").Value = 92

# Sheet 13: Chatham
$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("G3").Value = 7
$ws.Range("D7").Value = 12
$ws.Range("D9").Value = 50
$ws.Range("G9").Value = 35

# Sheet 14: Grand Crossing
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("C7").Value = 27
$ws.Range("C8").Value = 50
$ws.Range("C9").Value = 82

# Sheet 15: Loop
$ws = $wb.Worksheets.Item('Loop')
$ws.Range("F7").Value = 47
$ws.Range("C8").Value = 247
$ws.Range("E8").Value = 508
$ws.Range("F8").Value = 479
$ws.Range("C9").Value = 292
$ws.Range("E9").Value = 575
$ws.Range("F9").Value = 542

# Sheet 16: Armour Square
$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("F5").Value = 7
$ws.Range("F6").Value = 21
$ws.Range("I6").Value = 12
$ws.Range("F7").Value = 31
$ws.Range("I7").Value = 27

# Sheet 18: Little Italy, UIC
$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("F5").Value = 22
$ws.Range("F7").Value = 53

# Sheet 25: Rush & Division
$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("B5").Value = 29
$ws.Range("B6").Value = 34

# Sheet 26: Englewood
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("F7").Value = 31
$ws.Range("F9").Value = 94

# Sheet 27: Lake View
$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("C7").Value = 39
$ws.Range("H7").Value = 13
$ws.Range("C8").Value = 61
$ws.Range("H8").Value = 37

# Sheet 31: River North
$ws = $wb.Worksheets.Item('River North')
$ws.Range("H2").Value = 2
$ws.Range("H7").Value = 19

# Sheet 37: Woodlawn
$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("D6").Value = 14
$ws.Range("D7").Value = 21

# Sheet 40: Near South Side
$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("H3").Value = 1
$ws.Range("H6").Value = 12

# Sheet 41: Avalon Park
$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("G3").Value = 1
$ws.Range("G6").Value = 3

# Sheet 47: West Pullman
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 5

# Sheet 58: New City
$ws = $wb.Worksheets.Item('New City')
$ws.Range("I5").Value = 2
$ws.Range("I6").Value = 5

# Sheet 70: Belmont Cragin
$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("B5").Value = 4
$ws.Range("B7").Value = 9

# Sheet 76: West Elsdon
$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("E3").Value = 9
$ws.Range("E4").Value = 9
